# Scheduled runner update: refresh market-price-derived numeric columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) on several leve rows
# across the per-class sheets. Source data has no formulas - plain values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2229.392
$ws.Range("J17").Value = 2317.5305
$ws.Range("L17").Value = 6952.5915
$ws.Range("N17").Value = -7288.5915

$ws.Range("J94").Value = 3000
$ws.Range("L94").Value = 3000
$ws.Range("N94").Value = -3902

$ws.Range("H98").Value = 2788.3333
$ws.Range("I98").Value = 2687.111
$ws.Range("K98").Value = 2687.111
$ws.Range("M98").Value = -1189.111

$ws.Range("H122").Value = 2788.3333
$ws.Range("I122").Value = 2687.111
$ws.Range("K122").Value = 8061.333
$ws.Range("M122").Value = -5611.333

$ws.Range("H127").Value = 531.8570999999999
$ws.Range("I127").Value = 457.46155
$ws.Range("J127").Value = 1499
$ws.Range("K127").Value = 1372.38465
$ws.Range("L127").Value = 4497
$ws.Range("M127").Value = 3587.61535
$ws.Range("N127").Value = -14417

$ws.Range("H138").Value = 2694.05
$ws.Range("I138").Value = 2615.818
$ws.Range("J138").Value = 2789.6667
$ws.Range("K138").Value = 7847.454000000001
$ws.Range("L138").Value = 8369.000100000001
$ws.Range("M138").Value = -2707.454000000001
$ws.Range("N138").Value = -18649.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 54107.8
$ws.Range("I15").Value = 500
$ws.Range("J15").Value = 67509.75
$ws.Range("K15").Value = 500
$ws.Range("L15").Value = 67509.75
$ws.Range("M15").Value = -150
$ws.Range("N15").Value = -68209.75

$ws.Range("H61").Value = 2671.2778
$ws.Range("I61").Value = 1833.5652
$ws.Range("J61").Value = 4153.385
$ws.Range("K61").Value = 1833.5652
$ws.Range("L61").Value = 4153.385
$ws.Range("M61").Value = -1621.5652
$ws.Range("N61").Value = -4577.385

$ws.Range("H132").Value = 4355.4443
$ws.Range("I132").Value = 1768.5769
$ws.Range("J132").Value = 7895.3687
$ws.Range("K132").Value = 5305.7307
$ws.Range("L132").Value = 23686.1061
$ws.Range("M132").Value = -2775.7307
$ws.Range("N132").Value = -28746.1061

$ws.Range("H134").Value = 29400
$ws.Range("J134").Value = 29400
$ws.Range("L134").Value = 29400
$ws.Range("N134").Value = -39540

$ws.Range("H136").Value = 2671.2778
$ws.Range("I136").Value = 1833.5652
$ws.Range("J136").Value = 4153.385
$ws.Range("K136").Value = 5500.6956
$ws.Range("L136").Value = 12460.155
$ws.Range("M136").Value = -2950.6956
$ws.Range("N136").Value = -17560.155

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2812.9736
$ws.Range("I31").Value = 2035.762
$ws.Range("J31").Value = 3773.0588
$ws.Range("K31").Value = 2035.762
$ws.Range("L31").Value = 3773.0588
$ws.Range("M31").Value = -1740.762
$ws.Range("N31").Value = -4363.0588

$ws.Range("H34").Value = 2812.9736
$ws.Range("I34").Value = 2035.762
$ws.Range("J34").Value = 3773.0588
$ws.Range("K34").Value = 2035.762
$ws.Range("L34").Value = 3773.0588
$ws.Range("M34").Value = -1833.762
$ws.Range("N34").Value = -4177.0588

$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()

$ws.Range("H105").Value = 1578.9524
$ws.Range("I105").Value = 1342.9231
$ws.Range("J105").Value = 1962.5
$ws.Range("K105").Value = 1342.9231
$ws.Range("L105").Value = 1962.5
$ws.Range("M105").Value = 404.0769
$ws.Range("N105").Value = -5456.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 1789.4445
$ws.Range("I97").Value = 2950
$ws.Range("J97").Value = 1457.8572
$ws.Range("K97").Value = 8850
$ws.Range("L97").Value = 4373.571599999999
$ws.Range("M97").Value = -8354
$ws.Range("N97").Value = -5365.571599999999

$ws.Range("H98").Value = 1017.625
$ws.Range("J98").Value = 1378
$ws.Range("L98").Value = 4134
$ws.Range("N98").Value = -7130

$ws.Range("H107").Value = 411.75
$ws.Range("J107").Value = 421.63635
$ws.Range("L107").Value = 1264.90905
$ws.Range("N107").Value = -5104.90905

$ws.Range("H110").Value = 3357.1428
$ws.Range("I110").Value = 900
$ws.Range("J110").Value = 3766.6667
$ws.Range("K110").Value = 2700
$ws.Range("L110").Value = 11300.0001
$ws.Range("M110").Value = 1390
$ws.Range("N110").Value = -19480.0001

$ws.Range("H131").Value = 1065.6897
$ws.Range("I131").Value = 740
$ws.Range("J131").Value = 1103.2693
$ws.Range("K131").Value = 2220
$ws.Range("L131").Value = 3309.8079
$ws.Range("M131").Value = 2820
$ws.Range("N131").Value = -13389.8079

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 16352.25
$ws.Range("I70").Value = 10000
$ws.Range("J70").Value = 18469.666
$ws.Range("K70").Value = 10000
$ws.Range("L70").Value = 18469.666
$ws.Range("M70").Value = -9730
$ws.Range("N70").Value = -19009.666

$ws.Range("H73").Value = 16352.25
$ws.Range("I73").Value = 10000
$ws.Range("J73").Value = 18469.666
$ws.Range("K73").Value = 10000
$ws.Range("L73").Value = 18469.666
$ws.Range("M73").Value = -9064
$ws.Range("N73").Value = -20341.666

$ws.Range("H95").Value = 23344
$ws.Range("J95").Value = 23344
$ws.Range("L95").Value = 23344
$ws.Range("N95").Value = -28836

$ws.Range("H126").Value = 3057.081
$ws.Range("I126").Value = 3177.087
$ws.Range("J126").Value = 2859.9285
$ws.Range("K126").Value = 9531.261
$ws.Range("L126").Value = 8579.7855
$ws.Range("M126").Value = -7061.261
$ws.Range("N126").Value = -13519.7855

$ws.Range("H132").Value = 1986450.9
$ws.Range("I132").Value = 4168857.5
$ws.Range("J132").Value = 2444.7273
$ws.Range("K132").Value = 12506572.5
$ws.Range("L132").Value = 7334.1819
$ws.Range("M132").Value = -12504042.5
$ws.Range("N132").Value = -12394.1819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 753.86365
$ws.Range("J22").Value = 931
$ws.Range("L22").Value = 931
$ws.Range("N22").Value = -1521

$ws.Range("H24").Value = 2602
$ws.Range("I24").Value = 806
$ws.Range("J24").Value = 3500
$ws.Range("K24").Value = 806
$ws.Range("L24").Value = 3500
$ws.Range("M24").Value = -463
$ws.Range("N24").Value = -4186

$ws.Range("H27").Value = 753.86365
$ws.Range("J27").Value = 931
$ws.Range("L27").Value = 931
$ws.Range("N27").Value = -1145

